# Applies the "read me file added" diff:
#  - Sheet "Page_1": G1 header text change, G2:G13 descriptor text change,
#    plus a few Power/Lumens (D/E) cell value swaps.
#  - Sheet "Page_2": G1 header text change, G2:G13 descriptor text change.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Page_1")
$ws2 = $wb.Worksheets.Item("Page_2")

# ---- Page_1 ----
$ws1.Range("G1").Value = "table_header_position"

$newDesc1 = '"Part Number - Can be found on the top right position of the page"'
for ($r = 2; $r -le 13; $r++) {
    $ws1.Cells.Item($r, 7).Value = $newDesc1
}

# Power / Lumens column tweaks on Page_1
$ws1.Range("D3").Value  = "16W"
$ws1.Range("E3").Value  = "1600lm"

$ws1.Range("D11").Value = "16W"
$ws1.Range("E11").Value = "1600lm"

$ws1.Range("D12").Value = ""
$ws1.Range("E12").Value = ""

$ws1.Range("D13").Value = "16W"
$ws1.Range("E13").Value = "1600lm"

# ---- Page_2 ----
$ws2.Range("G1").Value = "table_header_position"

$newDesc2 = '"Part Number Description Dimensions Power Lumens Colour Temp. - Can be found on the right side of the page"'
for ($r = 2; $r -le 13; $r++) {
    $ws2.Cells.Item($r, 7).Value = $newDesc2
}
